# Agregados parte de la entrega: nuevas funcionalidades del sistema
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 13
$ws.Range("A13").Value = "AgregarDireccion"
$ws.Range("B13").Value = "UserController"
$ws.Range("C13").Value = 47

# Fila 14
$ws.Range("A14").Value = "ObtenerDireccion"
$ws.Range("B14").Value = "UserController"
$ws.Range("C14").Value = 17

# Fila 16 (columna B antes que A, tal como quedo registrado el orden de
# strings compartidos en el archivo original)
$ws.Range("B16").Value = "ItemController"
$ws.Range("A16").Value = "VerCompras"
$ws.Range("C16").Value = 24

# Fila 17
$ws.Range("A17").Value = "Comprar"
$ws.Range("B17").Value = "BuyOrderController"
$ws.Range("C17").Value = 19

# Fila 18
$ws.Range("A18").Value = "AgregarACarrito"
$ws.Range("B18").Value = "CartController"
$ws.Range("C18").Value = 27

# Fila 15 (agregada al final, por eso "verVentas" queda como el ultimo
# string nuevo agregado a la tabla de strings compartidos)
$ws.Range("A15").Value = "verVentas"
$ws.Range("B15").Value = "ItemController"
$ws.Range("C15").Value = 13

# Deja la seleccion activa en A16, como en el archivo final
$ws.Range("A16").Select()
